$wb = $excel.ActiveWorkbook

$wsSummary   = $wb.Worksheets.Item("Summary")
$wsRepayment = $wb.Worksheets.Item("Repayment schedule")

# --- Summary sheet: re-run of the disbursement/repayment totals ---
# F2 picks up F3's (General) number format since its new value (0) no longer
# needs the 2-decimal display.
$wsSummary.Range("F3").Copy() | Out-Null
$wsSummary.Range("F2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$wsSummary.Range("F2").Value = 0
$wsSummary.Range("A3").Value = 693.36
$wsSummary.Range("E3").Value = 693.36
$wsSummary.Range("F3").Value = 0

# --- Repayment schedule sheet: updated schedule rows ---
# Row 3 (installment #1)
$wsRepayment.Range("B3").Value = 90
$wsRepayment.Range("C3").Value = 42095
$wsRepayment.Range("F3").Value = 591.83000000000004
$wsRepayment.Range("G3").Value = 9408.17
$wsRepayment.Range("H3").Value = 295.89

# Row 4 (installment #2)
$wsRepayment.Range("B4").Value = 0
$wsRepayment.Range("C4").Value = 42095
$wsRepayment.Range("F4").Value = 887.72
$wsRepayment.Range("G4").Value = 8520.4500000000007
$wsRepayment.Range("H4").Value = 0

# Row 5 (installment #3)
$wsRepayment.Range("B5").Value = 0
$wsRepayment.Range("F5").Value = 887.72
$wsRepayment.Range("H5").Value = 0

# Row 12 (installment #10)
$wsRepayment.Range("B12").Value = 61
$wsRepayment.Range("C12").Value = 42339
$wsRepayment.Range("F12").Value = 834.74
$wsRepayment.Range("G12").Value = 1806.79
$wsRepayment.Range("H12").Value = 52.98

# Row 13 (installment #11)
$wsRepayment.Range("B13").Value = 0
$wsRepayment.Range("F13").Value = 887.72
$wsRepayment.Range("G13").Value = 919.07
$wsRepayment.Range("H13").Value = 0

# Row 14 (installment #12)
$wsRepayment.Range("F14").Value = 919.07
$wsRepayment.Range("H14").Value = 9.3699999999999992
$wsRepayment.Range("K14").Value = 928.44
$wsRepayment.Range("Q14").Value = 928.44

# --- Sheet activation / selection matching the saved view state ---
$wsRepayment.Activate()
$wsRepayment.Range("I15").Select() | Out-Null
